$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45:125 down to 46:126.
$ws.Rows("45:45").Insert()

# Populate the newly inserted row 45 with the new price-record data.
$ws.Range("A45").Value = 3
$ws.Range("B45").Value = "Femacal de La Calera"
$ws.Range("C45").Value = "Coquimbo"
$ws.Range("D45").Value = "2022-01-10"
$ws.Range("E45").Value = 5
$ws.Range("F45").Value = 100112030
$ws.Range("G45").Value = "Poroto granado"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 38
$ws.Range("K45").Value = 26000
$ws.Range("L45").Value = 26000
$ws.Range("M45").Value = 26000
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 1040
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
